# Working on sequence recorder
# Add a new translation entry (row 9) to the "Translation" sheet:
#   TEXT ID = SingleUseId6, TYPOGRAPHY NAME = Default,
#   ALIGNMENT = Center, DIRECTION = LTR, GB = <value>

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B9").Value = "SingleUseId6"
$ws.Range("C9").Value = "Default"
$ws.Range("D9").Value = "Center"
$ws.Range("E9").Value = "LTR"
$ws.Range("F9").Value = "<value>"
